# infoform-20106065.xlsx — content updates on the "INFO" sheet.
#
# 1. "Briefly describe your work on other team projects" answer (A11)
# 2. "Briefly describe any leadership..." answer (A15) gains a trailing period
# 3. Six schedule cells that said "Clases" now say "Ocupado"
# 4. Current selection moves to A13 (cosmetic, matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INFO")

$ws.Range("A11").Value = "Proyectos universitarios."
$ws.Range("A15").Value = "Ninguna."

$ws.Range("B36").Value = "Ocupado"
$ws.Range("D36").Value = "Ocupado"
$ws.Range("C42").Value = "Ocupado"
$ws.Range("E42").Value = "Ocupado"
$ws.Range("F46").Value = "Ocupado"
$ws.Range("D48").Value = "Ocupado"

$ws.Range("A13").Select()
